# ajustando para nao imprimir os valores None
# Append the missing product rows (74-80) to the price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("91697550", "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Elite Series A1 Convencional TCL", "3.699.00"),
    @("91697550", "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Elite Series A1 Convencional TCL", "3.699.00"),
    @("86839655", "Cabo Flexível 2,5mm 100m Azul 750V SIL Fios", "159.90"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00")
)

$startRow = 74
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]

    # Columns A and C hold numeric-looking text (product codes / prices) that
    # must stay text, matching the rest of the sheet where every cell is a
    # plain string. Force text formatting before writing so Excel doesn't
    # silently coerce values like "159.90" into the number 159.9.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]

    $ws.Cells.Item($r, 2).Value = $values[1]

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $values[2]
}
